# issue #5: add legislator_id, name, date into dataframe
#
# Adds three new columns to the "股票" (stock) worksheet -- the third
# worksheet in the workbook -- carrying the filing metadata that used to
# live only in the output file name / path:
#   H -> date            (constant "2011-11-17" for every data row)
#   I -> legislator_name (constant "謝國樑" for every data row)
#   J -> legislator_id   (constant 1387, numeric, for every data row)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

$lastRow = 14

# ---- header row (row 1) ----
$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

# Match the look of the existing header cells (bold font, thin box border,
# centered horizontally, aligned to the top) -- same formatting already
# used by columns B1:G1.
$header = $ws.Range("H1:J1")
$header.Font.Bold = $true
$header.Borders.LineStyle = 1
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160

# The date column holds a literal "yyyy-mm-dd" string, not a real date
# value, so force a text format on it first -- otherwise Excel's input
# parser would silently convert "2011-11-17" into a date serial number.
$dateRange = $ws.Range("H2:H" + $lastRow)
$dateRange.NumberFormat = "@"

# ---- data rows (rows 2-14) ----
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = "2011-11-17"
    $ws.Cells.Item($r, 9).Value = "謝國樑"
    $ws.Cells.Item($r, 10).Value = 1387
}

$excel.CutCopyMode = $false
